$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add a new "2023" column (T) by copying formatting from the 2022 column (S) ---
$ws.Cells.Item(4,19).Copy($ws.Cells.Item(4,20))
$ws.Cells.Item(5,19).Copy($ws.Cells.Item(5,20))
$ws.Cells.Item(6,19).Copy($ws.Cells.Item(6,20))
$ws.Cells.Item(7,19).Copy($ws.Cells.Item(7,20))

$ws.Cells.Item(4,20).Value = 2023
$ws.Cells.Item(5,20).Value = 50
$ws.Cells.Item(6,20).Value = 35.1
$ws.Cells.Item(7,20).Value = 21

# --- 2. Give the table its Kyrgyz title in A1 (matching the Russian/English titles in B1/C1) ---
$ws.Cells.Item(1,1).Value = "10.5.1.1 Финансылык туруктуулуктун көрсөткүчтөрү"
$ws.Cells.Item(1,2).WrapText = $true

# --- 3. Shrink the "Items"/footnote block (rows 8-10) down to 8pt ---
$ws.Range("A8:C10").Font.Size = 8

# Footnote superscript number + footnote text runs need to be touched per-run
# so the rich-text split inside the shared string is preserved.
$runs = @(
    @{r=9;  c=1; s1=1; l1=1; s2=2; l2=41},
    @{r=9;  c=2; s1=1; l1=2; s2=3; l2=26},
    @{r=9;  c=3; s1=1; l1=1; s2=2; l2=24},
    @{r=10; c=1; s1=1; l1=1; s2=2; l2=75},
    @{r=10; c=2; s1=1; l1=1; s2=2; l2=85},
    @{r=10; c=3; s1=1; l1=1; s2=2; l2=78}
)
foreach ($run in $runs) {
    $cell = $ws.Cells.Item($run.r, $run.c)
    $cell.Characters($run.s1, $run.l1).Font.Size = 8
    $cell.Characters($run.s2, $run.l2).Font.Size = 8
}

# --- 4. Row 10 gets an explicit height now that the text is smaller ---
$ws.Rows(10).RowHeight = 30

# --- 5. Columns A:C become one uniform width ---
$ws.Columns("A:C").ColumnWidth = 33.7

# --- 6. Reset the lingering selection left over from editing ---
$ws.Range("A1").Select()
